$d = $word.ActiveDocument

# Locate the paragraph that ends the "Автоматическое и дистанционное управление" bullet
# list (the one ending in "Утечки газа;") so we can add a new bullet right after it.
$sourcePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -match "Утечки газа") {
        $sourcePara = $p
        break
    }
}

# Insert a brand-new paragraph right after it; Word clones the paragraph/run
# formatting (pStyle "a3", numPr ilvl=1/numId=2, szCs 28) from the following
# paragraph context automatically.
$sourcePara.Range.InsertParagraphAfter() | Out-Null
$newPara = $sourcePara.Next()

$newPara.Range.Text = "Состоянии пробок (коррекция времени будильника)"

# Split the new text right after "Состоянии пробок " to drop the _GoBack
# bookmark in between the two sentences (moving it away from its previous
# location, at the very end of the document).
$splitRange = $newPara.Range.Duplicate
$splitRange.Find.Execute("Состоянии пробок ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitRange.Collapse(0) | Out-Null
$d.Bookmarks.Add("_GoBack", $splitRange) | Out-Null
